# Daily attendance processing - 2025-11-28 22:25:24
# Normalise the "Recorded By" (column G) audit list on the
# "Session Analysis Results" sheet: the recorder that logged the
# session first is moved to the back of the list (left-rotate by one)
# for every row whose G cell holds more than one recorder, except rows
# already in the canonical "admin@admin.com, System" order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if ($val -eq "") { continue }
    if ($val -eq "admin@admin.com, System") { continue }

    $parts = $val -split ", "
    if ($parts.Length -le 1) { continue }

    $rotated = ($parts[1..($parts.Length - 1)] + $parts[0]) -join ", "
    $cell.Value = $rotated
}
